# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# "JIMENA CAROLINA BALLESTEROS BALLESTA" (CC 1235048451) gains a second
# "Periodo Mora" (1902) alongside her existing 1901 period, and her two
# rows move up to sit right after DAYANNA's row (instead of after DAIRO's
# block). "DAIRO ALBERTO TRASLAVIÑA TORRES" (CC 1143397563) keeps his 10
# periods (2102-2111), now listed in ascending order, and his
# "Salario Basico" (column G) is corrected from 939249 to 877803.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New contents for rows 17-28 (B:Tipo Doc, C:N° Doc, D:Nombre, E:Periodo Mora,
# F:Valor Mora, G:Salario Basico). Row 16 (DAYANNA, 1901) is unchanged.
$rows = @(
    @{ Row = 17; Doc = "1235048451"; Nombre = "JIMENA CAROLINA BALLESTEROS BALLESTA"; Periodo = "1901"; Mora = 1104;  Salario = 877803 },
    @{ Row = 18; Doc = "1235048451"; Nombre = "JIMENA CAROLINA BALLESTEROS BALLESTA"; Periodo = "1902"; Mora = 17667; Salario = 877803 },
    @{ Row = 19; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2102"; Mora = 35112; Salario = 877803 },
    @{ Row = 20; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2103"; Mora = 35112; Salario = 877803 },
    @{ Row = 21; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2104"; Mora = 35112; Salario = 877803 },
    @{ Row = 22; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2105"; Mora = 35112; Salario = 877803 },
    @{ Row = 23; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2106"; Mora = 35112; Salario = 877803 },
    @{ Row = 24; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2107"; Mora = 35112; Salario = 877803 },
    @{ Row = 25; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2108"; Mora = 35112; Salario = 877803 },
    @{ Row = 26; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2109"; Mora = 35112; Salario = 877803 },
    @{ Row = 27; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2110"; Mora = 35112; Salario = 877803 },
    @{ Row = 28; Doc = "1143397563"; Nombre = "DAIRO ALBERTO TRASLAVIÑA TORRES";      Periodo = "2111"; Mora = 26919; Salario = 877803 }
)

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = "CC"
    $ws.Range("C" + $r.Row).Value = $r.Doc
    $ws.Range("D" + $r.Row).Value = $r.Nombre
    $ws.Range("E" + $r.Row).Value = $r.Periodo
    $ws.Range("F" + $r.Row).Value = $r.Mora
    $ws.Range("G" + $r.Row).Value = $r.Salario
}
